# The diff shows that row 23's and row 24's data got swapped, and
# likewise row 28's and row 29's data got swapped. Only columns
# A, B, D, E, F, G, H, Q, R differ between the two rows in each pair;
# the rest of the row (C, I-P, S onward) is identical between the pair,
# so swapping just those columns reproduces the target state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr1 = "$col" + "23"
    $addr2 = "$col" + "24"
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value2 = $v2
    $ws.Range($addr2).Value2 = $v1
}

foreach ($col in $cols) {
    $addr1 = "$col" + "28"
    $addr2 = "$col" + "29"
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value2 = $v2
    $ws.Range($addr2).Value2 = $v1
}
